$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Replace-LogoShapeName {
    param(
        [bool]$IsHeader,
        [int]$StoryIndex,
        [string]$Descr,
        [string]$NewName,
        [int]$Id,
        [int]$Cx,
        [int]$Cy
    )

    if ($IsHeader) {
        $hf = $sec.Headers.Item($StoryIndex)
    } else {
        $hf = $sec.Footers.Item($StoryIndex)
    }
    $rng = $hf.Range

    $ishapes = $rng.InlineShapes
    $shp = $ishapes.Item(1)
    $shpRange = $shp.Range

    # Collapse the shape's anchor range so the subsequent InsertXML
    # replaces it in place instead of inserting a duplicate alongside it.
    $shpRange.Text = ""

    $frag = '<?xml version="1.0"?>' +
        '<w:document ' +
        'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' +
        'xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" ' +
        'xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" ' +
        'xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" ' +
        'xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">' +
        '<w:body><w:p><w:r><w:drawing>' +
        '<wp:inline distB="0" distT="0" distL="0" distR="0">' +
        '<wp:extent cx="' + $Cx + '" cy="' + $Cy + '"/>' +
        '<wp:effectExtent b="0" l="0" r="0" t="0"/>' +
        '<wp:docPr descr="' + $Descr + '" id="' + $Id + '" name="' + $NewName + '"/>' +
        '<a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture">' +
        '<pic:pic><pic:nvPicPr>' +
        '<pic:cNvPr descr="' + $Descr + '" id="0" name="' + $NewName + '"/>' +
        '<pic:cNvPicPr preferRelativeResize="0"/>' +
        '</pic:nvPicPr>' +
        '<pic:blipFill><a:blip r:embed="rId1"/><a:srcRect b="0" l="0" r="0" t="0"/>' +
        '<a:stretch><a:fillRect/></a:stretch></pic:blipFill>' +
        '<pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="' + $Cx + '" cy="' + $Cy + '"/></a:xfrm>' +
        '<a:prstGeom prst="rect"/><a:ln/></pic:spPr>' +
        '</pic:pic></a:graphicData></a:graphic>' +
        '</wp:inline></w:drawing></w:r></w:p></w:body></w:document>'

    $shpRange.InsertXML($frag)
}

# footer1.xml (COM Footers.Item(2), id=2) - PearsonLogo, image2.png -> image1.png
Replace-LogoShapeName -IsHeader $false -StoryIndex 2 `
    -Descr 'Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png' `
    -NewName 'image1.png' -Id 2 -Cx 952500 -Cy 285750

# footer2.xml (COM Footers.Item(1), id=4) - PearsonLogo, image2.png -> image1.png
Replace-LogoShapeName -IsHeader $false -StoryIndex 1 `
    -Descr 'Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png' `
    -NewName 'image1.png' -Id 4 -Cx 952500 -Cy 285750

# header1.xml (COM Headers.Item(2), id=1) - BTec_Logo-Orange, image1.jpg -> image2.jpg
Replace-LogoShapeName -IsHeader $true -StoryIndex 2 `
    -Descr 'BTec_Logo-Orange' `
    -NewName 'image2.jpg' -Id 1 -Cx 914400 -Cy 277792

# header2.xml (COM Headers.Item(1), id=3) - BTec_Logo-Orange, image1.jpg -> image2.jpg
Replace-LogoShapeName -IsHeader $true -StoryIndex 1 `
    -Descr 'BTec_Logo-Orange' `
    -NewName 'image2.jpg' -Id 3 -Cx 914400 -Cy 277792

Write-Output "Logo shape names swapped."
